# BIS-768: Fixed XLS export test files
#
# The sample-type export/import test fixture gains a new "Unique" column
# (header in L4), mirroring the existing "Multivalued" column (K4) both in
# content and formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header cell with its text.
$ws.Range("L4").Value = "Unique"

# Give L4 the exact same formatting (bold header style) as the neighbouring
# K4 ("Multivalued") header cell.
$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)

# Reflect the author's final selection/cursor position in the sheet.
$ws.Range("L5").Select()
